$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Date in A1 changes (45753 -> 45755)
$ws.Range("A1").Value = 45755

# Row 2 (ORL vs ATL / Paolo Banchero)
$ws.Range("A2").Value = "ORL vs ATL"
$ws.Range("B2").Value = "Paolo Banchero - Sobre 12.5 Dobles intentados"
$ws.Range("C2").Value = "1,34x"
$ws.Range("D2").Value = "Cumplio en 9/10 https://prnt.sc/CFckrbnVt_jq"

# Row 3 (ORL vs ATL / Dyson Daniels)
$ws.Range("A3").Value = "ORL vs ATL"
$ws.Range("B3").Value = "Dyson Daniels - Sobre 3.5 Dobles realizados"
$ws.Range("C3").Value = "1,32x"
$ws.Range("D3").Value = "Cumplio en 8/10 (hizo 3) https://prnt.sc/43uAp3BLOrz0"
# Shorter text no longer needs the extra wrapped line, so the row shrinks
$ws.Rows.Item(3).RowHeight = 24.95

# Row 4 (ORL vs ATL / Zaccharie Risacher)
$ws.Range("A4").Value = "ORL vs ATL"
$ws.Range("B4").Value = "Zaccharie Risacher - Sobre 2.5 Dobles realizados"
$ws.Range("C4").Value = "1,57x"
$ws.Range("D4").Value = "Cumplio en 8/10 https://prnt.sc/DL3mo6VJdwF7"

# Row 5 (subtotal)
$ws.Range("C5").Value = "3,00x"

# Row 6 (PHX vs GSW / Ryan Dunn)
$ws.Range("A6").Value = "PHX vs GSW"
$ws.Range("B6").Value = "Ryan Dunn - Sobre 3.5 Triples intentados"
$ws.Range("C6").Value = "1,40x"
$ws.Range("D6").Value = "Cumplio en 9/10 https://prnt.sc/Ra03vgmlMulr"

# Row 7 (PHX vs GSW / Brandin Podziemski)
$ws.Range("A7").Value = "PHX vs GSW"
$ws.Range("B7").Value = "Brandin Podziemski - Sobre 4.5 Goles de campo realizados"
$ws.Range("C7").Value = "1,41x"
$ws.Range("D7").Value = "Cumplio en 8/10 (hizo 4) https://prnt.sc/NAoER9vq_erZ"

# Row 8 (PHX vs GSW / Tyus Jones)
$ws.Range("A8").Value = "PHX vs GSW"
$ws.Range("B8").Value = "Tyus Jones - Sobre 0.5 Dobles realizados"
$ws.Range("C8").Value = "1,41x"
$ws.Range("D8").Value = "Cumplio en 8/10 https://prnt.sc/jrNhgNL1MWeb"

# Row 9 (subtotal)
$ws.Range("C9").Value = "3,15x"

# Row 10 (IND vs WAS / Tyrese Haliburton)
$ws.Range("A10").Value = "IND vs WAS"
$ws.Range("B10").Value = "Tyrese Haliburton - Sobre 11.5 Gol de Campo intentado"
$ws.Range("C10").Value = "1,45x"
$ws.Range("D10").Value = "Cumplio en 9/10 https://prnt.sc/CWe8sHI8cfF2"

# Row 11 (IND vs WAS / Jordan Poole)
$ws.Range("A11").Value = "IND vs WAS"
$ws.Range("B11").Value = "Jordan Poole - Sobre 4.5 Goles de campo realizados"
$ws.Range("C11").Value = "1,40x"
$ws.Range("D11").Value = "Cumplio en 8/10 https://prnt.sc/BJNnffHSNgVB"

# Row 12 (IND vs WAS / Pascal Siakam)
$ws.Range("A12").Value = "IND vs WAS"
$ws.Range("B12").Value = "Pascal Siakam - Sobre 1.5 Tiros libres realizados"
$ws.Range("C12").Value = "1,43x"
$ws.Range("D12").Value = "Cumplio en 9/10 (en la unica no le hicieron faltas.) https://prnt.sc/VhKnmGxbc11_"

# Row 13 (subtotal)
$ws.Range("C13").Value = "3,70x"

# Row 14 (grand total / duplinha)
$ws.Range("C14").Value = "34,97x"
$ws.Range("D14").Value = "https://stake.com/sports/home?betId=69c67952-a153-4f8e-9fb8-7c0fd5211861&modal=bet"

# Update selection to match the saved view state
$ws.Range("G14").Select()
